# Apply edits described by the diff:
# - rename "Include ValueSets" -> "Include from Pure tone thresh"
# - rename "Include ValueSets 2" -> "Include from Diagnostic Audio"
# - update the Date metadata value
# - convert the two "Include ValueSets*" sheets from referencing a ValueSet URL
#   to referencing "All codes" from a CodeSystem (System URI row added)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsPure = $wb.Worksheets.Item("Include ValueSets")
$wsDiag = $wb.Worksheets.Item("Include ValueSets 2")

# 1. Update the Date value on the Metadata sheet (row 8, column B)
$wsMeta.Range("B8").Value = "2022-05-03T15:05:17-04:00"

# 2. Rename the two sheets
$wsPure.Name = "Include from Pure tone thresh"
$wsDiag.Name = "Include from Diagnostic Audio"

# 3. Update "Include from Pure tone thresh" sheet (was "Include ValueSets")
$wsPure.Range("A1").Value = "Codes"
$wsPure.Range("A2").Value = "All codes"

# blank separator row, keep the same border/style as the row above
$wsPure.Range("A2").Copy($wsPure.Range("A3"))
$wsPure.Range("A2").Copy($wsPure.Range("B3"))
$wsPure.Range("A3").Value = ""
$wsPure.Range("B3").Value = ""

# System URI row pointing at the new CodeSystem
$wsPure.Range("A2").Copy($wsPure.Range("A4"))
$wsPure.Range("A2").Copy($wsPure.Range("B4"))
$wsPure.Range("A4").Value = "System URI"
$wsPure.Range("B4").Value = "http://hl7.org/fhir/us/pacio-splasch/CodeSystem/SPLASCHPureToneThresholdAudiometryPanelCS"

# 4. Update "Include from Diagnostic Audio" sheet (was "Include ValueSets 2")
$wsDiag.Range("A1").Value = "Codes"
$wsDiag.Range("A2").Value = "All codes"

# blank separator row
$wsDiag.Range("A2").Copy($wsDiag.Range("A3"))
$wsDiag.Range("A2").Copy($wsDiag.Range("B3"))
$wsDiag.Range("A3").Value = ""
$wsDiag.Range("B3").Value = ""

# System URI row pointing at the new CodeSystem
$wsDiag.Range("A2").Copy($wsDiag.Range("A4"))
$wsDiag.Range("A2").Copy($wsDiag.Range("B4"))
$wsDiag.Range("A4").Value = "System URI"
$wsDiag.Range("B4").Value = "http://hl7.org/fhir/us/pacio-splasch/CodeSystem/SPLASCHDiagnosticAudiologyResultsPanelCS"

Write-Output "done"
